# Insert a new weekly price record for "Albahaca" (Primera, Región de Arica
# y Parinacota) as row 728 on Sheet1, shifting every existing row from 728
# downward (previously ending at 790) down by one (now ending at 791).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 728..790 down to 729..791 and leave a blank row 728 behind.
$ws.Rows.Item(728).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A728").Value = 6
$ws.Range("B728").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C728").Value = "Metropolitana"
$ws.Range("D728").Value = 45166
$ws.Range("E728").Value = 13
$ws.Range("F728").Value = 100112052
$ws.Range("G728").Value = "Albahaca"
$ws.Range("H728").Value = "Sin especificar"
$ws.Range("I728").Value = "Primera"
$ws.Range("J728").Value = 510
$ws.Range("K728").Value = 4500
$ws.Range("L728").Value = 5000
$ws.Range("M728").Value = 4755
$ws.Range("N728").Value = "$/paquete"
$ws.Range("O728").Value = "Región de Arica y Parinacota"
$ws.Range("P728").Value = 4755
$ws.Range("Q728").Value = 1
$ws.Range("R728").Value = "Hortaliza"
